$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New version entry details (row 22) - [1.17.8]
$detailsText = "reads the latest template only that contains AssignedTo`r`nall the other templates that don't have it in its place will not be read.`r`nthe extra columns that begin with ""A1-2"" will not be used in the generation since they are rare"

# Apply formatting to match the style used by the rest of the version log
# (A: left/top aligned, B: left/top aligned + wrap, C: left/top aligned date)
# *before* setting the values, so Excel doesn't auto-create a brand new
# number-format style for the date cell.
$ws.Range("A22:C22").HorizontalAlignment = -4131
$ws.Range("A22:C22").VerticalAlignment = -4160
$ws.Range("B22").WrapText = $true
$ws.Range("C22").NumberFormat = "d-mmm-yy"

$ws.Range("A22").Value = "[1.17.8]"
$ws.Range("B22").Value = $detailsText
$ws.Range("C22").Value = (Get-Date -Year 2018 -Month 9 -Day 21 -Hour 0 -Minute 0 -Second 0).Date

$ws.Rows.Item(22).RowHeight = 60

$ws.Range("C21").Select()
